# Update "săpt. 3" (column E) attendance values from 1 to 2 for the
# students listed below. The "Prezențe" (Q) column holds SUM(C:P) shared
# formulas, so it recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 6, 8, 12, 13, 14, 17, 18, 19, 21)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Value = 2
}

# Update the saved selection to match the author's last-edited cell.
$ws.Range("J13").Select()
